$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$win = $excel.ActiveWindow
$excel.Goto($ws.Range("A184"), $true)
Write-Host "ScrollRow:" $win.ScrollRow
Write-Host "ScrollColumn:" $win.ScrollColumn
